$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 63.6
$ws.Range("N2").Value = 54.82400714602223

$ws.Range("D3").Value = 107200
$ws.Range("E3").Value = 60.2
$ws.Range("F3").Value = 6.67
$ws.Range("K3").Value = 58.4
$ws.Range("N3").Value = 54.82400714602223

$ws.Range("K4").Value = 51.6
$ws.Range("N4").Value = 54.82400714602223

$ws.Range("D5").Value = 536000
$ws.Range("E5").Value = 32.1
$ws.Range("F5").Value = 1.13
$ws.Range("K5").Value = 48.8
$ws.Range("N5").Value = 54.82400714602223

$ws.Range("K6").Value = 40.8
$ws.Range("N6").Value = 54.82400714602223

$ws.Range("K7").Value = 39.6
$ws.Range("N7").Value = 54.82400714602223
